# Auto-generated Excel COM-interop script
# Applies targeted cell value updates to Sheet1 per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    "G2=110",
    "W2=1.01",
    "I3=1.26",
    "J3=8",
    "K3=9.4",
    "R3=1.81",
    "V3=4.8",
    "Z3=10",
    "AA3=10.5",
    "AD3=12",
    "AE3=16",
    "G4=2.92",
    "K4=3.2",
    "P4=1.67",
    "Q4=2.32",
    "T4=1.91",
    "AB4=9.6",
    "AD4=14",
    "AO4=44",
    "H5=7.8",
    "L5=1.26",
    "M5=1.03",
    "N5=5",
    "P5=2.36",
    "Q5=1.6",
    "R5=1.54",
    "T5=1.91",
    "U5=1.91",
    "Z5=100",
    "AA5=370",
    "AD5=40",
    "AE5=160",
    "AM5=150",
    "AO5=190",
    "F6=5.2",
    "G6=5.5",
    "H6=1.8",
    "I6=1.82",
    "J6=3.85",
    "K6=3.95",
    "O6=1.39",
    "S6=4",
    "T6=2",
    "U6=1.9",
    "V6=2.2",
    "W6=1.22",
    "AG6=21",
    "F7=1.77",
    "G7=1.88",
    "K7=4.2",
    "H8=1.78",
    "K8=4",
    "F9=4.8",
    "G9=6.2",
    "H9=1.62",
    "I9=1.78",
    "J9=4.2",
    "N9=5.2",
    "P9=2.44",
    "R9=1.58",
    "T9=1.62",
    "V9=2.28",
    "W9=1.2",
    "AO9=8.6",
    "I10=2.4",
    "V10=1.71",
    "AN10=60",
    "G11=10.5",
    "I11=1.4",
    "J11=5.5",
    "K11=5.6",
    "L11=1.38",
    "N11=4",
    "Q11=1.8",
    "V11=3.5",
    "X11=19",
    "Y11=8.199999999999999",
    "Z11=9",
    "AA11=12.5",
    "AB11=32",
    "AC11=11.5",
    "AD11=11.5",
    "AE11=18.5",
    "AF11=100",
    "AG11=40",
    "AH11=36",
    "AI11=50",
    "AJ11=490",
    "AK11=220",
    "AL11=180",
    "AM11=230",
    "AN11=320",
    "AO11=7.2",
    "L12=1.22",
    "F13=2.16",
    "S13=2.9",
    "G14=5",
    "L14=1.22",
    "M14=1.03",
    "Y14=15.5",
    "Z14=15.5",
    "AA14=22",
    "AN14=34",
    "M15=1.03",
    "Q15=1.45",
    "G16=4.2",
    "M16=1.02",
    "R16=2",
    "Y16=21",
    "AA16=25",
    "AC16=14.5",
    "AD16=13",
    "AE16=970",
    "AF16=40",
    "AG16=970",
    "AH16=970",
    "AI16=23",
    "AL16=38",
    "AM16=48",
    "AN16=970",
    "F17=3.3",
    "I17=2.52",
    "J17=3.35",
    "N17=3.2",
    "Q17=2.22",
    "R17=1.28",
    "U17=2",
    "V17=1.66",
    "W17=1.41",
    "Z17=15.5",
    "AA17=38",
    "AE17=48",
    "AH17=20",
    "AO17=29",
    "Y18=12",
    "Z18=28",
    "AB18=8.6",
    "AC18=7.6",
    "AD18=18.5",
    "AF18=16",
    "AG18=13.5",
    "AH18=25",
    "AK18=36",
    "AN18=36",
    "Y19=6.8",
    "AI19=70",
    "H20=2.02",
    "F21=1.79",
    "K21=3.75",
    "Q21=2.36",
    "R21=1.27",
    "W21=2.24",
    "AC21=8.4",
    "AD21=23",
    "I22=3.3",
    "J22=3.2",
    "K22=3.25",
    "N22=3.15",
    "U22=1.98",
    "V22=1.43",
    "W22=1.6",
    "F23=2.2",
    "O23=1.65",
    "R23=1.16",
    "AL23=80",
    "G24=16.5",
    "J24=6.6",
    "K24=6.8",
    "X24=19",
    "AB24=38",
    "AG24=60",
    "AO24=5.8",
    "G25=2.04",
    "I25=5.9",
    "L25=1.47",
    "P25=1.74",
    "Q25=2.18",
    "S25=3.6",
    "T25=1.97",
    "W25=1.98",
    "F27=3.6",
    "G27=3.8",
    "I27=2.58",
    "J27=2.96",
    "K27=3",
    "L27=1.67",
    "M27=1.16",
    "N27=2.22",
    "O27=1.75",
    "P27=1.39",
    "Q27=3.2",
    "R27=1.14",
    "S27=7.8",
    "T27=2.46",
    "U27=1.6",
    "W27=1.35",
    "X27=6.6",
    "Y27=6.4",
    "AB27=8.800000000000001",
    "AD27=14",
    "AE27=46",
    "AG27=19",
    "AH27=32",
    "AI27=95",
    "AJ27=95",
    "AK27=80",
    "AL27=130",
    "AM27=300",
    "AN27=130",
    "AO27=150",
    "H28=2.62",
    "I28=2.76",
    "K28=3",
    "I30=2.62",
    "L30=1.47",
    "V30=1.62",
    "AC30=7.6",
    "AC31=8.6",
    "Q34=2.34"
)

foreach ($u in $updates) {
    $parts = $u -split "="
    $addr = $parts[0]
    $val = [double]$parts[1]
    $ws.Range($addr).Value = $val
}
